# Updates the cryptos list sheet: refreshes price/volume values for Tue Nov 19 12:38:09 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "92.300.53"
Set-TextValue $ws "E2" "  +1.82%  "
Set-TextValue $ws "D3" "3.141.47"
Set-TextValue $ws "E3" "  +2.03%  "
Set-TextValue $ws "E4" "  +0.06%  "
Set-TextValue $ws "D5" "245.96"
Set-TextValue $ws "E5" "  +0.90%  "
Set-TextValue $ws "D6" "620.15"
Set-TextValue $ws "E6" "  +0.37%  "
Set-TextValue $ws "D7" "1.11"
Set-TextValue $ws "E7" "  -1.47%  "
Set-TextValue $ws "D8" "0.389"
Set-TextValue $ws "E8" "  +6.76%  "
Set-TextValue $ws "E9" "  -0.02%  "
Set-TextValue $ws "D10" "3.139.36"
Set-TextValue $ws "E10" "  +2.14%  "
Set-TextValue $ws "D11" "0.753"
Set-TextValue $ws "E11" "  +2.53%  "
Set-TextValue $ws "E12" "  +1.19%  "
Set-TextValue $ws "D13" "0.0000254"
Set-TextValue $ws "E13" "  +3.63%  "
Set-TextValue $ws "D14" "35.30"
Set-TextValue $ws "E14" "  +1.19%  "
Set-TextValue $ws "D15" "5.61"
Set-TextValue $ws "E15" "  +2.72%  "
Set-TextValue $ws "D16" "92.013.13"
Set-TextValue $ws "E16" "  +1.61%  "
Set-TextValue $ws "D17" "3.729.18"
Set-TextValue $ws "E17" "  +2.94%  "
Set-TextValue $ws "D18" "3.141.16"
Set-TextValue $ws "E18" "  +1.85%  "
Set-TextValue $ws "D19" "3.74"
Set-TextValue $ws "E19" "  +2.31%  "
Set-TextValue $ws "D20" "15.07"
Set-TextValue $ws "E20" "  +4.26%  "
Set-TextValue $ws "D21" "5.95"
Set-TextValue $ws "E21" "  +3.40%  "
Set-TextValue $ws "D22" "9.51"
Set-TextValue $ws "E22" "  +5.22%  "
Set-TextValue $ws "D23" "452.43"
Set-TextValue $ws "E23" "  +2.46%  "
Set-TextValue $ws "D24" "0.0000205"
Set-TextValue $ws "E24" "  -2.23%  "
Set-TextValue $ws "D25" "5.72"
Set-TextValue $ws "E25" "  +2.36%  "
Set-TextValue $ws "D26" "88.80"
Set-TextValue $ws "E26" "  -2.44%  "
Set-TextValue $ws "D27" "11.94"
Set-TextValue $ws "E27" "  +0.80%  "
Set-TextValue $ws "D28" "3.304.73"
Set-TextValue $ws "E28" "  +2.53%  "
Set-TextValue $ws "D29" "0.143"
Set-TextValue $ws "E29" "  +29.16%  "
Set-TextValue $ws "E30" "  +0.37%  "
Set-TextValue $ws "D31" "0.237"
Set-TextValue $ws "E31" "  +0.02%  "
Set-TextValue $ws "D32" "0.170"
Set-TextValue $ws "E32" "  -5.26%  "
Set-TextValue $ws "D33" "9.46"
Set-TextValue $ws "E33" "  +3.88%  "
Set-TextValue $ws "D34" "0.178"
Set-TextValue $ws "E34" "  +5.15%  "
Set-TextValue $ws "D35" "0.998"
Set-TextValue $ws "E35" "  +4.60%  "
Set-TextValue $ws "D36" "8.17"
Set-TextValue $ws "E36" "  +7.36%  "
Set-TextValue $ws "B37" "MantraDAO"
Set-TextValue $ws "C37" "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws "D37" "4.34"
Set-TextValue $ws "E37" "  +1.56%  "
Set-TextValue $ws "B38" "EthereumClassic"
Set-TextValue $ws "C38" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D38" "26.64"
Set-TextValue $ws "E38" "  +0.85%  "
Set-TextValue $ws "D39" "1.95"
Set-TextValue $ws "E39" "  +3.27%  "
Set-TextValue $ws "D40" "497.10"
Set-TextValue $ws "E40" "  +1.50%  "
Set-TextValue $ws "D41" "1.33"
Set-TextValue $ws "E41" "  +3.95%  "
Set-TextValue $ws "D42" "0.444"
Set-TextValue $ws "E42" "  +5.78%  "
Set-TextValue $ws "D43" "3.48"
Set-TextValue $ws "E43" "  -0.92%  "
Set-TextValue $ws "D44" "22.21"
Set-TextValue $ws "E44" "  +0.12%  "
Set-TextValue $ws "B46" "Stacks"
Set-TextValue $ws "C46" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D46" "1.94"
Set-TextValue $ws "E46" "  +3.03%  "
Set-TextValue $ws "B47" "Monero"
Set-TextValue $ws "C47" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D47" "157.89"
Set-TextValue $ws "E47" "  +2.61%  "
Set-TextValue $ws "B48" "ARBITRUM"
Set-TextValue $ws "C48" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D48" "0.710"
Set-TextValue $ws "E48" "  +4.17%  "
Set-TextValue $ws "D49" "1.38"
Set-TextValue $ws "E49" "  +4.00%  "
Set-TextValue $ws "B50" "VeChain"
Set-TextValue $ws "C50" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D50" "0.0335"
Set-TextValue $ws "E50" "  +7.35%  "
Set-TextValue $ws "B51" "Filecoin"
Set-TextValue $ws "C51" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D51" "4.46"
Set-TextValue $ws "E51" "  -0.11%  "
